# Update column F (dSF) values on the active sheet to match repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -5
$ws.Range("F3").Value = -10
$ws.Range("F4").Value = -5
$ws.Range("F6").Value = 2
$ws.Range("F8").Value = -3
$ws.Range("F11").Value = 0
